$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two source rows were removed from the data (missing-data resampling), shifting
# everything below them up: first the "RM 232" row (originally row 26), then the
# "SC 92" row (originally row 28, now row 27 after the first deletion).
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-apply the missing-data (imputation holdout) pattern to column E/F at the
# final, post-shift row positions.
$ws.Range("F5").ClearContents()
$ws.Range("F11").Value = 17.65
$ws.Range("E19").Value = -6.5
$ws.Range("F19").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("F25").Value = 16.6
$ws.Range("E27").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("E33").Value = -10.7
